# Refresh the coin-price table (Sheet1!A1:E51) to the latest scrape.
# Source workbook only ever stores Price/Volume/Coin/Link as literal text
# (never as native numbers/percentages), so every write below targets
# .Value as a string. Cells whose new text happens to *look* like a plain
# number (e.g. "0.9999") are pre-formatted as Text ("@") so Excel does not
# silently coerce them into a numeric cell and normalize away things like
# trailing zeros (e.g. "0.05330" -> 0.0533).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "26.512.71"
$ws.Range("E2").Value = "  +0.70%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.727.63"
$ws.Range("E3").Value = "  +0.48%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.97"
$ws.Range("E5").Value = "  +2.16%  "

# Row 6: USDC
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.05%  "

# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4809"
$ws.Range("E7").Value = "  +1.99%  "

# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2669"
$ws.Range("E8").Value = "  +1.64%  "

# Row 9: Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06214"
$ws.Range("E9").Value = "  -0.03%  "

# Row 10: WrappedEther
$ws.Range("D10").Value = "1.726.08"
$ws.Range("E10").Value = "  +0.41%  "

# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07156"
$ws.Range("E11").Value = "  +1.19%  "

# Row 12: Solana
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.67"
$ws.Range("E12").Value = "  +2.92%  "

# Row 13: Polygon
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6179"
$ws.Range("E13").Value = "  +4.50%  "

# Row 14: Polkadot
$ws.Range("E14").Value = "  +2.47%  "

# Row 15: Litecoin
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.13"
$ws.Range("E15").Value = "  +0.99%  "

# Row 16: Dai
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9998"
$ws.Range("E16").Value = "  -0.05%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "26.524.35"
$ws.Range("E17").Value = "  +0.72%  "

# Row 18: BinanceUSD
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9997"
$ws.Range("E18").Value = "  -0.09%  "

# Row 19: ShibaInu
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006924"
$ws.Range("E19").Value = "  +1.91%  "

# Row 20: Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.65"
$ws.Range("E20").Value = "  +0.41%  "

# Row 21: WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "1.947.92"
$ws.Range("E21").Value = "  +0.78%  "

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.532"
$ws.Range("E22").Value = "  -0.57%  "

# Row 23: Cosmos
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.958"
$ws.Range("E23").Value = "  +1.91%  "

# Row 24: Chainlink
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.279"
$ws.Range("E24").Value = "  -1.09%  "

# Row 25: Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.17"
$ws.Range("E25").Value = "  +0.96%  "

# Row 26: EthereumClassic
$ws.Range("E26").Value = "  +0.99%  "

# Row 27: LidoDAOToken
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.798"
$ws.Range("E27").Value = "  +2.03%  "

# Row 28: Toncoin
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.403"
$ws.Range("E28").Value = "  -0.09%  "

# Row 29: BitcoinCash
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.73"
$ws.Range("E29").Value = "  -0.06%  "

# Row 30: InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.983"
$ws.Range("E30").Value = "  -1.09%  "

# Row 31: Stellar
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08021"
$ws.Range("E31").Value = "  +3.79%  "

# Row 32: Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.708"
$ws.Range("E32").Value = "  +0.42%  "

# Row 33: Hedera
$ws.Range("E33").Value = "  +2.70%  "

# Row 34: HuobiToken
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.614"
$ws.Range("E34").Value = "  +0.03%  "

# Row 35: ImmutableX
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6359"
$ws.Range("E35").Value = "  +2.58%  "

# Row 36: ARBITRUM
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9908"
$ws.Range("E36").Value = "  +1.83%  "

# Row 37: TrustWalletToken
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9329"
$ws.Range("E37").Value = "  +0.43%  "

# Row 38: RenderToken
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.099"
$ws.Range("E38").Value = "  +10.24%  "

# Row 39: MXToken
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.420"
$ws.Range("E39").Value = "  +0.33%  "

# Row 40: Quant
$ws.Range("B40").Value = "Quant"
$ws.Range("C40").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "105.10"
$ws.Range("E40").Value = "  -9.09%  "

# Row 41: PaxDollar
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.003"
$ws.Range("E41").Value = "  +0.19%  "

# Row 42: VeChain
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01501"
$ws.Range("E42").Value = "  +2.11%  "

# Row 43: FraxShare
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.638"
$ws.Range("E43").Value = "  +6.45%  "

# Row 44: TheSandbox
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3909"
$ws.Range("E44").Value = "  +2.43%  "

# Row 45: Aptos
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.919"
$ws.Range("E45").Value = "  +10.55%  "

# Row 46: Algorand
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1187"
$ws.Range("E46").Value = "  +2.63%  "

# Row 47: Cronos
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05330"
$ws.Range("E47").Value = "  +0.71%  "

# Row 48: Elrond
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.92"
$ws.Range("E48").Value = "  +0.97%  "

# Row 49: EnergySwap
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.871"
$ws.Range("E49").Value = "  +3.02%  "

# Row 50: NEARProtocol
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.267"
$ws.Range("E50").Value = "  +3.88%  "

# Row 51: Decentraland
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3429"
$ws.Range("E51").Value = "  +1.15%  "
